# Update reviews for parisk
# Split existing row 12 into two rows: clear C12 (polite_expressions) and
# insert a new row 13 containing the original row-12 data but with a new
# id/source_file/text/review_result (and C13 keeps "nan").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate row 12 into a brand-new row 13 (copy + insert shifts nothing
# below it, since row 13 is currently the first blank row).
$ws.Rows.Item(12).Copy()
$ws.Rows.Item(13).Insert()

# Row 12: clear the polite_expressions (C) cell.
$ws.Range("C12").Value = ""

# Row 13: update id, source_file, text, review_result to the new review entry.
$ws.Range("F13").Value = "77474e59-42ef-43e4-850b-a07d6b41a266"
$ws.Range("G13").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H13").Value = "You absolutely know this but you hide these results."
$ws.Range("I13").Value = "Correct"
